$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 13890463
$ws.Range("I62").Value = 18520042
$ws.Range("K62").Value = 18520042
$ws.Range("M62").Value = -18519418
$ws.Range("H65").Value = 13890463
$ws.Range("I65").Value = 18520042
$ws.Range("K65").Value = 92600210
$ws.Range("M65").Value = -92597090
$ws.Range("H69").Value = 5498.5293
$ws.Range("I69").Value = 4925
$ws.Range("J69").Value = 5675
$ws.Range("K69").Value = 14775
$ws.Range("L69").Value = 17025
$ws.Range("M69").Value = -13901
$ws.Range("N69").Value = -18773
$ws.Range("H72").Value = 5498.5293
$ws.Range("I72").Value = 4925
$ws.Range("J72").Value = 5675
$ws.Range("K72").Value = 44325
$ws.Range("L72").Value = 51075
$ws.Range("M72").Value = -39957
$ws.Range("N72").Value = -59811
$ws.Range("H74").Value = 5000
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -4064
$ws.Range("H77").Value = 5000
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 25000
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -20320
$ws.Range("H80").Value = 428.6316
$ws.Range("J80").Value = 475
$ws.Range("L80").Value = 1425
$ws.Range("N80").Value = -3421
$ws.Range("H83").Value = 428.6316
$ws.Range("J83").Value = 475
$ws.Range("L83").Value = 4275
$ws.Range("N83").Value = -14259
$ws.Range("H129").Value = 969.9545000000001
$ws.Range("J129").Value = 980.38336
$ws.Range("L129").Value = 2941.15008
$ws.Range("N129").Value = -12941.15008
$ws.Range("H132").Value = 1165.9661
$ws.Range("I132").Value = 903.9048
$ws.Range("J132").Value = 1813.4117
$ws.Range("K132").Value = 2711.7144
$ws.Range("L132").Value = 5440.2351
$ws.Range("M132").Value = -181.7143999999998
$ws.Range("N132").Value = -10500.2351
$ws.Range("H135").Value = 1400.807
$ws.Range("I135").Value = 1079.234
$ws.Range("J135").Value = 2912.2
$ws.Range("K135").Value = 9713.106
$ws.Range("L135").Value = 26209.8
$ws.Range("M135").Value = -7178.106
$ws.Range("N135").Value = -31279.8
$ws.Range("H141").Value = 1488.5306
$ws.Range("I141").Value = 990.4737
$ws.Range("J141").Value = 3209.0908
$ws.Range("K141").Value = 2971.4211
$ws.Range("L141").Value = 9627.2724
$ws.Range("M141").Value = 2208.5789
$ws.Range("N141").Value = -19987.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1286.5581
$ws.Range("I61").Value = 1178.5625
$ws.Range("J61").Value = 1600.7273
$ws.Range("K61").Value = 1178.5625
$ws.Range("L61").Value = 1600.7273
$ws.Range("M61").Value = -966.5625
$ws.Range("N61").Value = -2024.7273
$ws.Range("H74").Value = 1259.7273
$ws.Range("I74").Value = 954
$ws.Range("J74").Value = 2635.5
$ws.Range("K74").Value = 954
$ws.Range("L74").Value = 2635.5
$ws.Range("M74").Value = -80
$ws.Range("N74").Value = -4383.5
$ws.Range("H77").Value = 1259.7273
$ws.Range("I77").Value = 954
$ws.Range("J77").Value = 2635.5
$ws.Range("K77").Value = 4770
$ws.Range("L77").Value = 13177.5
$ws.Range("M77").Value = -402
$ws.Range("N77").Value = -21913.5
$ws.Range("H110").Value = 1023.75
$ws.Range("I110").Value = 931.6667
$ws.Range("J110").Value = 1300
$ws.Range("K110").Value = 931.6667
$ws.Range("L110").Value = 1300
$ws.Range("M110").Value = 1113.3333
$ws.Range("H132").Value = 2225054
$ws.Range("I132").Value = 1818.0883
$ws.Range("K132").Value = 5454.2649
$ws.Range("M132").Value = -2924.2649
$ws.Range("H136").Value = 1286.5581
$ws.Range("I136").Value = 1178.5625
$ws.Range("J136").Value = 1600.7273
$ws.Range("K136").Value = 3535.6875
$ws.Range("L136").Value = 4802.1819
$ws.Range("M136").Value = -985.6875
$ws.Range("N136").Value = -9902.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 12653.479
$ws.Range("I20").Value = 1789.6428
$ws.Range("J20").Value = 29552.777
$ws.Range("K20").Value = 1789.6428
$ws.Range("L20").Value = 29552.777
$ws.Range("M20").Value = -1542.6428
$ws.Range("N20").Value = -30046.777
$ws.Range("H80").Value = 287.5238
$ws.Range("I80").Value = 101.5
$ws.Range("J80").Value = 331.29413
$ws.Range("K80").Value = 101.5
$ws.Range("L80").Value = 331.29413
$ws.Range("M80").Value = 896.5
$ws.Range("N80").Value = -2327.29413
$ws.Range("H83").Value = 287.5238
$ws.Range("I83").Value = 101.5
$ws.Range("J83").Value = 331.29413
$ws.Range("K83").Value = 507.5
$ws.Range("L83").Value = 1656.47065
$ws.Range("M83").Value = 4484.5
$ws.Range("N83").Value = -11640.47065
$ws.Range("H99").Value = 45455684
$ws.Range("I99").Value = 58824390
$ws.Range("K99").Value = 58824390
$ws.Range("M99").Value = -58822892
$ws.Range("H107").Value = 100001270
$ws.Range("I107").Value = 200001180
$ws.Range("J107").Value = 1364.6
$ws.Range("K107").Value = 200001180
$ws.Range("L107").Value = 1364.6
$ws.Range("M107").Value = -199999260
$ws.Range("N107").Value = -5204.6
$ws.Range("H134").Value = 2057.432
$ws.Range("I134").Value = 1668.6
$ws.Range("J134").Value = 2569.0527
$ws.Range("K134").Value = 5005.799999999999
$ws.Range("L134").Value = 7707.158100000001
$ws.Range("M134").Value = -2470.799999999999
$ws.Range("N134").Value = -12777.1581

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6987.26
$ws.Range("I31").Value = 1767.9395
$ws.Range("J31").Value = 17118.883
$ws.Range("K31").Value = 1767.9395
$ws.Range("L31").Value = 17118.883
$ws.Range("M31").Value = -1472.9395
$ws.Range("N31").Value = -17708.883
$ws.Range("H34").Value = 6987.26
$ws.Range("I34").Value = 1767.9395
$ws.Range("J34").Value = 17118.883
$ws.Range("K34").Value = 1767.9395
$ws.Range("L34").Value = 17118.883
$ws.Range("M34").Value = -1565.9395
$ws.Range("N34").Value = -17522.883
$ws.Range("H58").Value = 1249.017
$ws.Range("I58").Value = 645.46344
$ws.Range("J58").Value = 2623.7778
$ws.Range("K58").Value = 645.46344
$ws.Range("L58").Value = 2623.7778
$ws.Range("M58").Value = -442.46344
$ws.Range("N58").Value = -3029.7778
$ws.Range("H107").Value = 15152118
$ws.Range("I107").Value = 18519140
$ws.Range("J107").Value = 517.5
$ws.Range("K107").Value = 18519140
$ws.Range("L107").Value = 517.5
$ws.Range("M107").Value = -18517220
$ws.Range("N107").Value = -4357.5
$ws.Range("H132").Value = 2069.78
$ws.Range("I132").Value = 1623.75
$ws.Range("J132").Value = 2862.7222
$ws.Range("K132").Value = 4871.25
$ws.Range("L132").Value = 8588.1666
$ws.Range("M132").Value = -2341.25
$ws.Range("N132").Value = -13648.1666
$ws.Range("H136").Value = 1249.017
$ws.Range("I136").Value = 645.46344
$ws.Range("J136").Value = 2623.7778
$ws.Range("K136").Value = 1936.39032
$ws.Range("L136").Value = 7871.3334
$ws.Range("M136").Value = 613.60968
$ws.Range("N136").Value = -12971.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3547.0625
$ws.Range("I122").Value = 466.66666
$ws.Range("J122").Value = 4752.4346
$ws.Range("K122").Value = 4199.99994
$ws.Range("L122").Value = 42771.9114
$ws.Range("M122").Value = -1749.99994
$ws.Range("N122").Value = -47671.9114

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6620.3125
$ws.Range("J70").Value = 5754.5
$ws.Range("L70").Value = 5754.5
$ws.Range("N70").Value = -6294.5
$ws.Range("H73").Value = 6620.3125
$ws.Range("J73").Value = 5754.5
$ws.Range("L73").Value = 5754.5
$ws.Range("N73").Value = -7626.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 102230.4
$ws.Range("I7").Value = 113089.336
$ws.Range("K7").Value = 113089.336
$ws.Range("M7").Value = -112977.336
$ws.Range("H126").Value = 102230.4
$ws.Range("I126").Value = 113089.336
$ws.Range("K126").Value = 339268.008
$ws.Range("M126").Value = -336798.008
$ws.Range("H132").Value = 22267418
$ws.Range("I132").Value = 31434366
$ws.Range("J132").Value = 4828.2856
$ws.Range("K132").Value = 94303098
$ws.Range("L132").Value = 14484.8568
$ws.Range("M132").Value = -94300568
$ws.Range("N132").Value = -19544.8568
$ws.Range("H136").Value = 4799.364
$ws.Range("I136").Value = 1985.9767
$ws.Range("J136").Value = 14880.667
$ws.Range("K136").Value = 5957.9301
$ws.Range("L136").Value = 44642.001
$ws.Range("M136").Value = -3407.9301
$ws.Range("N136").Value = -49742.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 61177556
$ws.Range("I107").Value = 100001400
$ws.Range("J107").Value = 5714921
$ws.Range("K107").Value = 300004200
$ws.Range("L107").Value = 17144763
$ws.Range("M107").Value = -300002280
$ws.Range("N107").Value = -17148603
$ws.Range("H136").Value = 6412299
$ws.Range("I136").Value = 1916.2982
$ws.Range("J136").Value = 23811910
$ws.Range("K136").Value = 5748.8946
$ws.Range("L136").Value = 71435730
$ws.Range("M136").Value = -3198.8946
$ws.Range("N136").Value = -71440830
